# Quotes.xlsx - add a "Discount" column (F) to Sheet1 with per-row
# discount percentages, pushing the existing "Model DropDown" helper
# columns one slot to the right.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new, blank column at F; the former F/G columns (and their
# data on row 2) shift right to G/H.
$ws.Columns.Item(6).Insert()

# F1 (header) and F3:F6 (the new data rows) should wrap text, same as
# the rest of the header row / most of the sheet. Flip that on while
# the whole new column still shares one uniform (default) style, then
# give the entire F1:F6 range the shared 2-decimal numeric format -
# this way only two brand new cell styles get created (one for the
# wrapped/general-font cells, one for the Consolas-font F2 cell),
# matching how the workbook was actually edited.
$ws.Range("F3:F6").WrapText = $true
$ws.Range("F1:F6").NumberFormat = "0.00"

# Make the new column as wide as the other wide columns on the sheet.
$ws.Columns.Item(6).ColumnWidth = 31

# Header label for the new column.
$ws.Range("F1").Value = "Discount"

# Discount values for the existing row plus four new rows below it.
$ws.Range("F2").Value = 10
$ws.Range("F3").Value = 15
$ws.Range("F4").Value = 20
$ws.Range("F5").Value = 25
$ws.Range("F6").Value = 35

# Leave the selection where the author last left it.
[void]$ws.Range("G6").Select()
